# Auto-generated update of cryptos list (price/volume refresh + 3 row pair swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.215.39'
$ws.Range('E2').Value = '  +10.63%  '
$ws.Range('D3').Value = '3.483.03'
$ws.Range('E3').Value = '  +6.96%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '416.67'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.81%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '123.48'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +13.74%  '
$ws.Range('D7').Value = '3.473.50'
$ws.Range('E7').Value = '  +6.87%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.600'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.88%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.684'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +10.64%  '
$ws.Range('E11').Value = '  +40.07%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '41.62'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +6.11%  '
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '4.037.50'
$ws.Range('E14').Value = '  +6.98%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.61'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +4.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.08'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +6.02%  '
$ws.Range('D17').Value = '3.484.49'
$ws.Range('E17').Value = '  +7.75%  '
$ws.Range('D18').Value = '63.196.47'
$ws.Range('E18').Value = '  +10.89%  '
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.05'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0000144'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +34.37%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.31'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '317.74'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +7.96%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.43'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +11.61%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '31.20'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +11.36%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.78'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.36%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.92'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('B30').Value = 'LEO'
$ws.Range('C30').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.31'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.57%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.174'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +3.80%  '
$ws.Range('E32').Value = '  +3.17%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.63'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +4.17%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '42.11'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.62%  '
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.55'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +19.36%  '
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.01'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.68%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0489'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.30%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '52.05'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.49'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('E42').Value = '  +7.55%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.126'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.92%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '135.65'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '16.96'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.64%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.282'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.90'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.27'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.66%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '22.00'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('D50').Value = '3.827.85'
$ws.Range('E50').Value = '  +6.88%  '
$ws.Range('D51').Value = '2.188.03'
$ws.Range('E51').Value = '  +1.99%  '
